# Fix duplicate data and correct reconciliation results.
#
# The sheet had 5 columns (A:TradeID, B:FundingCurve, C:DiscountBasis,
# D:CSA_Type, E:ModelVersion) and 25 data rows, but rows 7-16 (T006-T015)
# had been duplicated into rows 17-26. This script replaces the data with
# the corrected, de-duplicated 25-trade table (T001-T025), drops the
# DiscountBasis column entirely, and renames/reorders the remaining
# columns to TradeID, FundingCurve, CSA_Type, ModelVersion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The final table has only 4 columns (A:D) -- remove the old column C
# (DiscountBasis) so the old D/E (CSA_Type/ModelVersion) data shift left
# into C/D.
$ws.Columns.Item(3).Delete()

# Corrected, de-duplicated data: header + 25 unique trades (T001-T025).
$data = @(
    @("TradeID", "FundingCurve", "CSA_Type", "ModelVersion"),
    @("T001", "USD-LIBOR", "Cleared", "v2024.1"),
    @("T002", "EUR-LIBOR", "Bilateral", "v2024.2"),
    @("T003", "USD-LIBOR", "Cleared", "v2024.1"),
    @("T004", "EUR-LIBOR", "Bilateral", "v2024.2"),
    @("T005", "USD-LIBOR", "Cleared", "v2024.1"),
    @("T006", "USD-LIBOR", "Cleared", "v2024.1"),
    @("T007", "EUR-LIBOR", "Bilateral", "v2024.2"),
    @("T008", "USD-LIBOR", "Cleared", "v2024.1"),
    @("T009", "USD-LIBOR", "Cleared", "v2024.3"),
    @("T010", "EUR-LIBOR", "Bilateral", "v2024.2"),
    @("T011", "USD-LIBOR", "Cleared", "v2024.1"),
    @("T012", "EUR-LIBOR", "Bilateral", "v2024.2"),
    @("T013", "USD-LIBOR", "Cleared", "v2024.3"),
    @("T014", "EUR-LIBOR", "Bilateral", "v2024.1"),
    @("T015", "USD-LIBOR", "Cleared", "v2024.2"),
    @("T016", "EUR-LIBOR", "Bilateral", "v2024.1"),
    @("T017", "USD-LIBOR", "Cleared", "v2024.2"),
    @("T018", "EUR-LIBOR", "Bilateral", "v2024.1"),
    @("T019", "USD-LIBOR", "Cleared", "v2024.3"),
    @("T020", "EUR-LIBOR", "Bilateral", "v2024.2"),
    @("T021", "USD-LIBOR", "Cleared", "v2024.1"),
    @("T022", "EUR-LIBOR", "Bilateral", "v2024.2"),
    @("T023", "USD-LIBOR", "Cleared", "v2024.3"),
    @("T024", "EUR-LIBOR", "Bilateral", "v2024.1"),
    @("T025", "USD-LIBOR", "Cleared", "v2024.2")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}

# Remove any leftover content beyond the new table bounds (in case of
# stale values past column D or row 26 from the old duplicated data).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count
if ($lastCol -gt 4) {
    $ws.Range($ws.Cells.Item(1, 5), $ws.Cells.Item($lastRow, $lastCol)).Clear()
}
if ($lastRow -gt 26) {
    $ws.Range($ws.Cells.Item(27, 1), $ws.Cells.Item($lastRow, 4)).Clear()
}
